# Update dashboards - 2025-10-30
# Applies the data refresh for rows 29, 30, 48, 49, 50, 52 on the
# "Aguilar Prototype" sheet: bump the "as of" dates and refresh the
# trailing-values history (columns Q:U) to reflect the new pull.
#
# The "Latest Date" cells (column N) are stored as literal text (e.g.
# "2025-10-29"), not real dates, so they are written with a leading
# apostrophe to force text entry and keep Excel from re-interpreting
# them as date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: 5yr, 5yr Forward (T5YIFR)
$ws.Range("N29").Value = "'2025-10-29"
$ws.Range("Q29").Value = 2.23

# Row 30: 10yr TIPS (T10YIE)
$ws.Range("N30").Value = "'2025-10-29"
$ws.Range("Q30").Value = 2.3
$ws.Range("S30").Value = 2.28
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.3

# Row 48: 2y UST (DGS2)
$ws.Range("N48").Value = "'2025-10-28"
$ws.Range("Q48").Value = 3.47
$ws.Range("T48").Value = 3.48

# Row 49: 5y UST (DGS5)
$ws.Range("N49").Value = "'2025-10-28"
$ws.Range("Q49").Value = 3.6
$ws.Range("T49").Value = 3.61

# Row 50: 10y UST (DGS10)
$ws.Range("N50").Value = "'2025-10-28"
$ws.Range("Q50").Value = 3.99
$ws.Range("R50").Value = 4.01
$ws.Range("S50").Value = 4.02
$ws.Range("T50").Value = 4.01
$ws.Range("U50").Value = 3.97

# Row 52: BAA (DBAA)
$ws.Range("N52").Value = "'2025-10-28"
$ws.Range("R52").Value = 5.64
$ws.Range("T52").Value = 5.67
$ws.Range("U52").Value = 5.66
